{"js": "// Database release note update: populate the first empty data row of the\n// \"Database Changes\" table with the new SQL release script's details, and\n// rename the script per the new naming convention.\n//\n// Layout of the table (row index 0 = header row):\n//   rows.items[1] (first data row)  -> Serial Number | Script Name | Run in\n//                                      Database | Note\n//   rows.items[2] (second data row) -> left blank, but this is where the\n//                                      cursor ends up after the edit, so\n//                                      Word's \"_GoBack\" bookmark (tracking\n//                                      the last edit position) moves from\n//                                      the end of the document into this\n//                                      row's \"Script Name\" cell.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 2 (index 1): fill in the release note for the new drop-orphan-tables\n// script.\nconst dataRow = rows.items[1];\nconst dataCells = dataRow.cells;\ndataCells.load(\"items\");\nawait context.sync();\n\ndataCells.items[0].body.insertText(\"001\", \"End\");\ndataCells.items[1].body.insertText(\"001_09062015DropTables.sql\", \"End\");\ndataCells.items[2].body.insertText(\"aidr_fetch_manager\", \"End\");\ndataCells.items[3].body.insertText(\"Drop orphan tables\", \"End\");\n\n// Row 3 (index 2), \"Script Name\" cell (index 1): this is where Word leaves\n// the \"_GoBack\" bookmark after the edit above. There is only ever one\n// \"_GoBack\" bookmark in a document, so remove the old one (at the end of\n// the document, after the two soft hyphens) before adding the new one.\nconst nextRow = rows.items[2];\nconst nextRowCells = nextRow.cells;\nnextRowCells.load(\"items\");\nawait context.sync();\n\ncontext.document.deleteBookmark(\"_GoBack\");\nnextRowCells.items[1].body.getRange(\"End\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Database release note update: populate the first empty data row of the\n# \"Database Changes\" table with the new SQL release script's details, and\n# rename the script per the new naming convention.\n#\n# Layout of the table (row 1 = header):\n#   Row 2 (first data row)  -> Serial Number | Script Name | Run in Database | Note\n#   Row 3 (second data row) -> left untouched (still blank) except that this\n#                               is where the cursor/selection ends up after\n#                               the edit, so Word's \"_GoBack\" bookmark (last\n#                               edit position) moves from the end of the\n#                               document into this row's \"Script Name\" cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(2, 1).Range.InsertAfter(\"001\")\n$t.Cell(2, 2).Range.InsertAfter(\"001_09062015DropTables.sql\")\n$t.Cell(2, 3).Range.InsertAfter(\"aidr_fetch_manager\")\n$t.Cell(2, 4).Range.InsertAfter(\"Drop orphan tables\")\n\n# Word keeps a single \"_GoBack\" bookmark tracking the last edit location;\n# re-adding it here removes the previous one (at the end of the document)\n# and places it at the new last-edit spot.\n$d.Bookmarks.Add(\"_GoBack\", $t.Cell(3, 2).Range)\n"}
